# Updates the Internacional.xlsx round-26 historical sheet:
#  - Column A (match id) gets new numeric ids
#  - Column E ("Round") switches from text "Matchweek N" to plain numeric N

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ A = 3;  E = 2 }
    3  = @{ A = 4;  E = 3 }
    4  = @{ A = 7;  E = 5 }
    5  = @{ A = 11; E = 8 }
    6  = @{ A = 14; E = 10 }
    7  = @{ A = 18; E = 13 }
    8  = @{ A = 20; E = 15 }
    9  = @{ A = 22; E = 17 }
    10 = @{ A = 24; E = 18 }
    11 = @{ A = 27; E = 20 }
    12 = @{ A = 32; E = 23 }
    13 = @{ A = 35; E = 25 }
    14 = @{ A = 37; E = 26 }
    15 = @{ A = 28; E = 19 }
    16 = @{ E = 16 }
    17 = @{ A = 6;  E = 7 }
    18 = @{ A = 28; E = 21 }
    19 = @{ A = 5;  E = 1 }
    20 = @{ A = 19; E = 14 }
    21 = @{ A = 31; E = 24 }
    22 = @{ A = 12; E = 6 }
    23 = @{ A = 6;  E = 4 }
    24 = @{ A = 12; E = 9 }
    25 = @{ A = 29; E = 22 }
    26 = @{ A = 10; E = 11 }
    27 = @{ A = 16; E = 12 }
}

foreach ($rowNum in $updates.Keys) {
    $vals = $updates[$rowNum]
    if ($vals.ContainsKey('A')) {
        $ws.Range("A$rowNum").Value = $vals['A']
    }
    if ($vals.ContainsKey('E')) {
        $ws.Range("E$rowNum").Value = $vals['E']
    }
}
